$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values could be auto-detected as numbers by Excel.
# Force them to remain text (matching the original inlineStr cell type)
# by temporarily switching the cell to Text format, then clearing the
# formatting again afterwards so no extra style is left on the cell.
$textForcedValues = [ordered]@{
    "D5" = "592.15"
    "D6" = "136.20"
    "D11" = "5.30"
    "D14" = "33.98"
    "D20" = "471.51"
    "D21" = "14.09"
    "D23" = "7.70"
    "D28" = "7.89"
    "D29" = "6.94"
    "D32" = "26.56"
    "D34" = "2.51"
    "D36" = "5.77"
    "D37" = "52.10"
    "D39" = "0.0385"
    "D40" = "419.95"
    "D45" = "0.265"
    "D48" = "25.36"
    "D49" = "0.113"
    "D51" = "119.89"
}

foreach ($addr in $textForcedValues.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textForcedValues[$addr]
    $cell.ClearFormats()
}

# Remaining cells: the new text is not number-like, so a plain assignment
# already keeps Excel from re-interpreting it as a number/date.
$plainValues = [ordered]@{
    "D2" = "63.052.37"
    "E2" = "  -2.17%  "
    "D3" = "3.118.22"
    "E3" = "  -0.74%  "
    "E4" = "  +0.04%  "
    "E5" = "  -2.84%  "
    "E6" = "  -5.42%  "
    "E7" = "  +0.09%  "
    "D8" = "3.112.63"
    "E8" = "  -0.83%  "
    "E9" = "  -2.04%  "
    "E10" = "  -4.08%  "
    "E11" = "  -1.66%  "
    "E12" = "  -3.07%  "
    "E13" = "  -4.93%  "
    "E14" = "  -4.02%  "
    "D15" = "3.631.60"
    "E15" = "  -0.67%  "
    "E16" = "  +2.22%  "
    "D17" = "62.983.91"
    "E17" = "  -2.19%  "
    "D18" = "3.123.45"
    "E18" = "  -0.53%  "
    "E19" = "  -3.11%  "
    "E20" = "  -1.42%  "
    "E21" = "  -5.00%  "
    "E22" = "  -3.25%  "
    "E23" = "  -0.78%  "
    "E24" = "  +0.65%  "
    "E25" = "  -4.26%  "
    "E26" = "  +0.04%  "
    "E27" = "  -1.98%  "
    "E28" = "  -6.56%  "
    "E29" = "  -4.08%  "
    "E31" = "  +0.04%  "
    "E32" = "  -1.00%  "
    "E33" = "  -7.22%  "
    "E34" = "  -4.75%  "
    "E35" = "  -2.84%  "
    "E36" = "  -3.56%  "
    "E37" = "  -0.81%  "
    "E38" = "  -8.90%  "
    "B39" = "VeChain"
    "C39" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "E39" = "  -1.97%  "
    "B40" = "Bittensor"
    "C40" = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
    "E40" = "  -6.06%  "
    "E41" = "  -1.06%  "
    "D42" = "2.890.87"
    "E42" = "  +0.12%  "
    "E43" = "  -11.92%  "
    "E44" = "  -6.37%  "
    "E45" = "  +1.07%  "
    "E47" = "  -5.98%  "
    "E48" = "  -3.44%  "
    "E49" = "  -0.53%  "
    "E50" = "  -5.85%  "
    "E51" = "  +0.16%  "
}

foreach ($addr in $plainValues.Keys) {
    $ws.Range($addr).Value = $plainValues[$addr]
}

